# Update "Horarios Linea 141" workbook with the latest scrape (03:53:17).

$wb = $excel.ActiveWorkbook

# --- Sheet 1: LP1912 -----------------------------------------------------
$ws1 = $wb.Worksheets.Item("LP1912")

$ws1.Range("A2").Value = "Última actualización: 03:53:17"
$ws1.Range("A3").Value = "Total filas: 7"

# Row 6
$ws1.Range("A6").Value = "03:53:17"
$ws1.Range("B6").Value = "04:02"
$ws1.Range("C6").Value = "81_EL PELIGRO"
$ws1.Range("D6").Value = 9
$ws1.Range("E6").Value = "LP1912"

# Row 7
$ws1.Range("A7").Value = "03:53:17"
$ws1.Range("B7").Value = "04:48"
$ws1.Range("C7").Value = "81_EL PELIGRO"
$ws1.Range("D7").Value = 55
$ws1.Range("E7").Value = "LP1912"

# Row 8
$ws1.Range("A8").Value = "03:53:17"
$ws1.Range("B8").Value = "04:53"
$ws1.Range("C8").Value = "11_ETCHEVERRY"
$ws1.Range("D8").Value = 60
$ws1.Range("E8").Value = "LP1912"

# Row 9
$ws1.Range("A9").Value = "03:53:17"
$ws1.Range("B9").Value = "05:17"
$ws1.Range("C9").Value = "17_ROMERO"
$ws1.Range("D9").Value = 84
$ws1.Range("E9").Value = "LP1912"

# Row 10
$ws1.Range("A10").Value = "03:53:17"
$ws1.Range("B10").Value = "05:22"
$ws1.Range("C10").Value = "23_HERNANDEZ"
$ws1.Range("D10").Value = 89
$ws1.Range("E10").Value = "LP1912"

# Row 11 (new)
$ws1.Range("A11").Value = "03:53:17"
$ws1.Range("B11").Value = "05:44"
$ws1.Range("C11").Value = "14_ABASTO"
$ws1.Range("D11").Value = 111
$ws1.Range("E11").Value = "LP1912"

# Row 12 (new)
$ws1.Range("A12").Value = "03:53:17"
$ws1.Range("B12").Value = "05:47"
$ws1.Range("C12").Value = "17_ROMERO"
$ws1.Range("D12").Value = 114
$ws1.Range("E12").Value = "LP1912"

# --- Sheet 2: LP1912-215 --------------------------------------------------
$ws2 = $wb.Worksheets.Item("LP1912-215")
$ws2.Range("A2").Value = "Última actualización: 03:53:17"

# --- Sheet 3: 6203-6173 ----------------------------------------------------
$ws3 = $wb.Worksheets.Item("6203-6173")
$ws3.Range("A2").Value = "Última actualización: 03:53:17"
